# Scheduled runner update: refresh market-board derived price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) on each
# crafting-job sheet. Values are plain numbers (no formulas in the
# workbook), so each changed cell is written directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 759.4
$ws.Range("I28").Value = 749.3333
$ws.Range("J28").Value = 850
$ws.Range("K28").Value = 749.3333
$ws.Range("L28").Value = 850
$ws.Range("M28").Value = -264.3333
$ws.Range("N28").Value = -1820

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H58").Value = 505.9
$ws.Range("I58").Value = 505.9
$ws.Range("K58").Value = 1517.7
$ws.Range("M58").Value = -1367.7

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 859.125
$ws.Range("I61").Value = 859.125
$ws.Range("K61").Value = 2577.375
$ws.Range("M61").Value = -2405.375

$ws.Range("H88").Value = 16222.4
$ws.Range("J88").Value = 19403
$ws.Range("L88").Value = 19403
$ws.Range("N88").Value = -20215

$ws.Range("H91").Value = 16222.4
$ws.Range("J91").Value = 19403
$ws.Range("L91").Value = 19403
$ws.Range("N91").Value = -22211

$ws.Range("H100").Value = 4238.5625
$ws.Range("I100").Value = 2608.077
$ws.Range("K100").Value = 2608.077
$ws.Range("M100").Value = -2067.077

$ws.Range("H138").Value = 2901.6667
$ws.Range("J138").Value = 2984.125
$ws.Range("L138").Value = 8952.375
$ws.Range("N138").Value = -19232.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7289.3706
$ws.Range("I2").Value = 1629.2
$ws.Range("J2").Value = 14364.583
$ws.Range("K2").Value = 1629.2
$ws.Range("L2").Value = 14364.583
$ws.Range("M2").Value = -1516.2
$ws.Range("N2").Value = -14590.583

$ws.Range("H60").Value = 19232.334
$ws.Range("J60").Value = 19173
$ws.Range("L60").Value = 19173
$ws.Range("N60").Value = -20639

$ws.Range("H61").Value = 1591618.1
$ws.Range("I61").Value = 2384410.2
$ws.Range("K61").Value = 2384410.2
$ws.Range("M61").Value = -2384198.2

$ws.Range("H88").Value = 1640
$ws.Range("I88").Value = 1975.25
$ws.Range("J88").Value = 1472.375
$ws.Range("K88").Value = 1975.25
$ws.Range("L88").Value = 1472.375
$ws.Range("M88").Value = -1569.25
$ws.Range("N88").Value = -2284.375

$ws.Range("H91").Value = 1640
$ws.Range("I91").Value = 1975.25
$ws.Range("J91").Value = 1472.375
$ws.Range("K91").Value = 1975.25
$ws.Range("L91").Value = 1472.375
$ws.Range("M91").Value = -571.25
$ws.Range("N91").Value = -4280.375

$ws.Range("H116").Value = 7289.3706
$ws.Range("I116").Value = 1629.2
$ws.Range("J116").Value = 14364.583
$ws.Range("K116").Value = 1629.2
$ws.Range("L116").Value = 14364.583
$ws.Range("M116").Value = 664.8
$ws.Range("N116").Value = -18952.583

$ws.Range("H136").Value = 1591618.1
$ws.Range("I136").Value = 2384410.2
$ws.Range("K136").Value = 7153230.600000001
$ws.Range("M136").Value = -7150680.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7289.3706
$ws.Range("I3").Value = 1629.2
$ws.Range("J3").Value = 14364.583
$ws.Range("K3").Value = 1629.2
$ws.Range("L3").Value = 14364.583
$ws.Range("M3").Value = -1515.2
$ws.Range("N3").Value = -14592.583

$ws.Range("H20").Value = 2047.6364
$ws.Range("I20").Value = 1275.4
$ws.Range("K20").Value = 1275.4
$ws.Range("M20").Value = -1028.4

$ws.Range("H86").Value = 2700
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377

$ws.Range("H87").Value = 102489
$ws.Range("I87").Value = 79978
$ws.Range("K87").Value = 79978
$ws.Range("M87").Value = -78730

$ws.Range("H88").Value = 17133.2
$ws.Range("J88").Value = 18338.75
$ws.Range("L88").Value = 18338.75
$ws.Range("N88").Value = -19150.75

$ws.Range("H89").Value = 2700
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884

$ws.Range("H90").Value = 102489
$ws.Range("I90").Value = 79978
$ws.Range("K90").Value = 239934
$ws.Range("M90").Value = -233694

$ws.Range("H91").Value = 17133.2
$ws.Range("J91").Value = 18338.75
$ws.Range("L91").Value = 18338.75
$ws.Range("N91").Value = -21146.75

$ws.Range("H94").Value = 809.2174
$ws.Range("I94").Value = 631.7778
$ws.Range("J94").Value = 1448
$ws.Range("K94").Value = 631.7778
$ws.Range("L94").Value = 1448
$ws.Range("M94").Value = -180.7778
$ws.Range("N94").Value = -2350

$ws.Range("H99").Value = 7796.3184
$ws.Range("I99").Value = 3707
$ws.Range("J99").Value = 14952.625
$ws.Range("K99").Value = 3707
$ws.Range("L99").Value = 14952.625
$ws.Range("M99").Value = -2209
$ws.Range("N99").Value = -17948.625

$ws.Range("H107").Value = 6853.923
$ws.Range("I107").Value = 7091.8335
$ws.Range("K107").Value = 7091.8335
$ws.Range("M107").Value = -5171.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 231.21053
$ws.Range("I7").Value = 119.07692
$ws.Range("J7").Value = 474.16666
$ws.Range("K7").Value = 119.07692
$ws.Range("L7").Value = 474.16666
$ws.Range("M7").Value = -6.076920000000001
$ws.Range("N7").Value = -700.16666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 13602.333
$ws.Range("I87").Value = 9451.909
$ws.Range("K87").Value = 28355.727
$ws.Range("M87").Value = -27107.727

$ws.Range("H90").Value = 13602.333
$ws.Range("I90").Value = 9451.909
$ws.Range("K90").Value = 85067.181
$ws.Range("M90").Value = -78827.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 37400.4
$ws.Range("J33").Value = 35500.5
$ws.Range("L33").Value = 35500.5
$ws.Range("N33").Value = -36004.5

$ws.Range("H40").Value = 21800.8
$ws.Range("J40").Value = 21800.8
$ws.Range("L40").Value = 21800.8
$ws.Range("N40").Value = -22102.8

$ws.Range("H43").Value = 7508.778
$ws.Range("J43").Value = 18177.334
$ws.Range("L43").Value = 18177.334
$ws.Range("N43").Value = -18479.334

$ws.Range("H94").Value = 48614.332
$ws.Range("J94").Value = 48614.332
$ws.Range("L94").Value = 48614.332
$ws.Range("N94").Value = -49966.332

$ws.Range("H113").Value = 3014.8572
$ws.Range("I113").Value = 1815.8334
$ws.Range("J113").Value = 3914.125
$ws.Range("K113").Value = 1815.8334
$ws.Range("L113").Value = 3914.125
$ws.Range("M113").Value = 354.1666
$ws.Range("N113").Value = -8254.125

$ws.Range("H122").Value = 47308.082
$ws.Range("I122").Value = 79549.84
$ws.Range("J122").Value = 9204.182000000001
$ws.Range("K122").Value = 238649.52
$ws.Range("L122").Value = 27612.546
$ws.Range("M122").Value = -236199.52
$ws.Range("N122").Value = -32512.546

$ws.Range("H123").Value = 45000.5
$ws.Range("J123").Value = 45000.5
$ws.Range("L123").Value = 45000.5
$ws.Range("N123").Value = -49900.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2023.7858
$ws.Range("I46").Value = 1056
$ws.Range("J46").Value = 2185.0833
$ws.Range("K46").Value = 1056
$ws.Range("L46").Value = 2185.0833
$ws.Range("M46").Value = -868
$ws.Range("N46").Value = -2561.0833

$ws.Range("H132").Value = 2183161.5
$ws.Range("I132").Value = 6964737.5
$ws.Range("K132").Value = 20894212.5
$ws.Range("M132").Value = -20891682.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 1000
$ws.Range("M39").Value = -587

$ws.Range("H122").Value = 3719.074
$ws.Range("I122").Value = 2592.1052
$ws.Range("J122").Value = 6395.625
$ws.Range("K122").Value = 7776.3156
$ws.Range("L122").Value = 19186.875
$ws.Range("M122").Value = -5326.3156
$ws.Range("N122").Value = -24086.875

$ws.Range("H132").Value = 11189777
$ws.Range("I132").Value = 25157440
$ws.Range("K132").Value = 75472320
$ws.Range("M132").Value = -75469790

